# Added reset test result function
# Resets all recorded CRUD test results on the "Test Results" sheet back to FALSE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")
$ws.Activate()

# Reset every recorded test result (Create/Read/Update/Delete Test Passed) to FALSE.
$ws.Range("B2:E24").Value = $false

# Move the active selection to where the user left off after running the reset.
$ws.Range("H14").Select()
